$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 37: DOH -> DOH, 7, DESERTDAYS27, 1
$ws.Range("A37").Value = "DOH"
$ws.Range("B37").Value = "DOH"
$ws.Range("C37").Value = 7
$ws.Range("D37").Value = "DESERTDAYS27"
$ws.Range("E37").Value = 1

# Row 38: DUB -> DUB, 7, DESERTDAYSDUBAI27, 1
$ws.Range("A38").Value = "DUB"
$ws.Range("B38").Value = "DUB"
$ws.Range("C38").Value = 7
$ws.Range("D38").Value = "DESERTDAYSDUBAI27"
$ws.Range("E38").Value = 1

# Update selection to match the new active cell location
$ws.Activate()
$ws.Range("D39").Select()
